$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("C2").Value = 2.120147338132171
$ws.Range("D2").Value = 0.032664625996545
$ws.Range("E2").Value = 2.786581942191909
$ws.Range("F2").Value = 0.08492224995459739

# Row 3
$ws.Range("E3").Value = 2.575046244071891
$ws.Range("F3").Value = 0.07758764563809811

# Row 4
$ws.Range("C4").Value = 0.4267577744766148
$ws.Range("D4").Value = 0.0103734718657173
$ws.Range("E4").Value = 0.5420737743758589
$ws.Range("F4").Value = 0.01492484181948727

# Row 5
$ws.Range("C5").Value = 0.2348088910634413
$ws.Range("D5").Value = 0.01015904091871061
$ws.Range("E5").Value = 0.2607220301076594
$ws.Range("F5").Value = 0.01311748503567509

# Row 6
$ws.Range("E6").Value = 0.2202390871830772
$ws.Range("F6").Value = 0.01215053414235768

# Row 7
$ws.Range("E7").Value = 0.2037080907976631
$ws.Range("F7").Value = 0.01124374412555548

# Row 8
$ws.Range("E8").Value = 0.1535832490935778
$ws.Range("F8").Value = 0.007908761114158546
